$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue {
    param($addr, $val)
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue "D2" "97.448.22"
Set-TextValue "E2" "  +0.81%  "
Set-TextValue "D3" "3.723.97"
Set-TextValue "E3" "  +0.12%  "
Set-TextValue "E4" "  -0.17%  "
Set-TextValue "D5" "2.17"
Set-TextValue "E5" "  +12.78%  "
Set-TextValue "D6" "237.67"
Set-TextValue "E6" "  -0.54%  "
Set-TextValue "D7" "656.08"
Set-TextValue "E7" "  +0.06%  "
Set-TextValue "D8" "0.442"
Set-TextValue "E8" "  +3.96%  "
Set-TextValue "E9" "  +3.79%  "
Set-TextValue "E10" "  -0.07%  "
Set-TextValue "D11" "3.722.54"
Set-TextValue "E11" "  +0.19%  "
Set-TextValue "D12" "0.0000312"
Set-TextValue "E12" "  +15.85%  "
Set-TextValue "D13" "44.70"
Set-TextValue "E13" "  -1.46%  "
Set-TextValue "E14" "  +0.53%  "
Set-TextValue "D15" "6.83"
Set-TextValue "E15" "  -0.29%  "
Set-TextValue "D16" "4.417.79"
Set-TextValue "E16" "  +0.12%  "
Set-TextValue "D17" "97.284.01"
Set-TextValue "E17" "  +0.87%  "
Set-TextValue "D18" "9.23"
Set-TextValue "E18" "  +1.61%  "
Set-TextValue "D19" "3.722.46"
Set-TextValue "E19" "  -0.25%  "
Set-TextValue "E20" "  +1.72%  "
Set-TextValue "D21" "18.87"
Set-TextValue "E21" "  -1.19%  "
Set-TextValue "D22" "0.536"
Set-TextValue "E22" "  +1.23%  "
Set-TextValue "D23" "528.17"
Set-TextValue "E23" "  +1.09%  "
Set-TextValue "D24" "3.46"
Set-TextValue "E24" "  -1.04%  "
Set-TextValue "D25" "0.0000224"
Set-TextValue "E25" "  +9.48%  "
Set-TextValue "D26" "117.24"
Set-TextValue "E26" "  +14.19%  "
Set-TextValue "E27" "  -2.24%  "
Set-TextValue "D28" "0.213"
Set-TextValue "E28" "  +26.28%  "
Set-TextValue "D29" "13.44"
Set-TextValue "E29" "  +0.41%  "
Set-TextValue "D30" "12.72"
Set-TextValue "E30" "  +1.08%  "
Set-TextValue "E31" "  -1.33%  "
Set-TextValue "E32" "  -0.18%  "
Set-TextValue "D33" "0.190"
Set-TextValue "E33" "  +2.53%  "
Set-TextValue "B34" "Fetch.AI"
Set-TextValue "C34" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D34" "1.83"
Set-TextValue "E34" "  -3.05%  "
Set-TextValue "B35" "EthereumClassic"
Set-TextValue "C35" "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D35" "33.11"
Set-TextValue "E35" "  +0.91%  "
Set-TextValue "D36" "0.999"
Set-TextValue "E36" "  -0.19%  "
Set-TextValue "D37" "0.597"
Set-TextValue "E37" "  -1.23%  "
Set-TextValue "D38" "638.55"
Set-TextValue "E38" "  -3.17%  "
Set-TextValue "D39" "8.73"
Set-TextValue "E39" "  -2.18%  "
Set-TextValue "E40" "  +0.01%  "
Set-TextValue "E41" "  +3.41%  "
Set-TextValue "B42" "Algorand"
Set-TextValue "C42" "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D42" "0.498"
Set-TextValue "E42" "  +9.33%  "
Set-TextValue "B43" "Filecoin"
Set-TextValue "C43" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D43" "6.82"
Set-TextValue "E43" "  -4.40%  "
Set-TextValue "B44" "EnergySwap"
Set-TextValue "C44" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D44" "40.87"
Set-TextValue "E44" "  +1.31%  "
Set-TextValue "D45" "2.00"
Set-TextValue "E45" "  +0.98%  "
Set-TextValue "D46" "0.966"
Set-TextValue "E46" "  -1.19%  "
Set-TextValue "D47" "0.0455"
Set-TextValue "E47" "  -0.49%  "
Set-TextValue "D48" "2.38"
Set-TextValue "E48" "  +1.94%  "
Set-TextValue "B49" "Cosmos"
Set-TextValue "C49" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue "D49" "8.77"
Set-TextValue "E49" "  +1.92%  "
Set-TextValue "B50" "WhiteBITCoin"
Set-TextValue "C50" "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue "D50" "23.64"
Set-TextValue "E50" "  +0.08%  "
Set-TextValue "E51" "  +3.71%  "
